# Adds a "note to Goodfellow" continuation to the document:
#  1. The hidden "_GoBack" bookmark that sat alone in the final (empty)
#     paragraph is moved up to sit right after the title run.
#  2. That now-empty final paragraph is expanded into several new
#     paragraphs of notes, ending with a "Historical trends in deep
#     learning" heading and a "P11" paragraph.

$d = $word.ActiveDocument

# --- Part 1: relocate the _GoBack bookmark onto the title paragraph ---
$titlePara = $d.Paragraphs(1)
$titleXml = @'
<w:p w14:paraId="4CBAAB5C" w14:textId="4AFD542B" w:rsidR="00CE345F" w:rsidRPr="00926C80" w:rsidRDefault="00A239FB" w:rsidP="00A239FB"><w:pPr><w:pStyle w:val="Rubrik"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r w:rsidRPr="00926C80"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Ian Goodfellow &#8211; Deep Learning</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@
$titlePara.Range.InsertXML($titleXml) | Out-Null

# --- Part 2: turn the trailing bookmark-only paragraph into the notes ---
$notesXml = @'
<w:p w14:paraId="6641C7C7" w14:textId="77777777" w:rsidR="00926C80" w:rsidRPr="00926C80" w:rsidRDefault="00926C80" w:rsidP="00A239FB"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">When using </w:t></w:r><w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>representation learning</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> you can detect the representation itself, hence not needing any predefined features. Usually renders better performance than hand designed representations. The quintessential example of a representation learning algorithm is the </w:t></w:r><w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>autoencoder</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>. They are trained to keep as much info as possible but also make the new representation have various nice properties.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">When designing features and algorithms we want to separate the </w:t></w:r><w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>factors of variation</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>The quintessential example of a deep learning model is the feedforward deep</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">network or </w:t></w:r><w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>multilayer perceptron (MLP)</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> A multilayer perceptron is just a</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>mathematical function mapping some set of input values to output values.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Deep learning is a </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>particular kind</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> of machine learning that achieves great power and flexibility by learning to represent the world as a nested hierarchy of concepts, with each concept defined in relation to simpler concepts, and more abstract representations computed in terms of less abstract ones</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="Rubrik1"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:lastRenderedPageBreak/><w:t>Historical trends in deep learning</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>P11</w:t></w:r></w:p>
'@
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertXML($notesXml) | Out-Null
